$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "39.910.19"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -3.72%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.330.81"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -4.67%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "308.45"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -2.96%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "84.80"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -6.33%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.527"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -2.97%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.483"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -3.68%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -2.22%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "29.96"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -7.34%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.72%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.692.46"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -4.44%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.42"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -4.85%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.67"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -4.68%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.336.82"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -4.19%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.756"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.91%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "39.877.40"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -3.46%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0901"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -2.65%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.10"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -3.23%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "67.63"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -6.21%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.60"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -6.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.48"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.54"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -6.16%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.15%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.80"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -4.93%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "23.20"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -4.12%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.12"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -5.10%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.22"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -3.98%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "33.83"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -3.69%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "153.60"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.98%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.03%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.10"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -3.98%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.43"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -4.30%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0716"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -4.32%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.114"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.02%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0996"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.28%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.73"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -6.95%  "
$ws.Range("B39").Value = "Celestia"
$ws.Range("C39").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "15.52"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -7.35%  "
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.71"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -4.45%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -2.93%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.942.59"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.29%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -4.84%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "17.65"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -3.28%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0262"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -5.26%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.24"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -2.70%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.70"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -6.64%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.553.19"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -4.68%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "92.44"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -3.56%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "70.37"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -4.74%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "63.08"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -3.35%  "
